$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 44, pushing the existing rows 44-61 down to 45-62.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with this week's new price entry.
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44466
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100106
$ws.Range("H44").Value = "Oleaginosos"
$ws.Range("I44").Value = 100106002
$ws.Range("J44").Value = "Palta"
$ws.Range("K44").Value = "Hass"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 400
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 26000
$ws.Range("P44").Value = 25500
$ws.Range("Q44").Value = "$/bandeja 10 kilos"
$ws.Range("R44").Value = "Perú"
$ws.Range("S44").Value = 2550
$ws.Range("T44").Value = 10
